$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "Toewijzingen"

# --- Row 2: add new RFID / 3V3 labels in H2 / I2 (no special fill) ---
$ws.Range("H2").Value = "RFID"
$ws.Range("H2").HorizontalAlignment = -4108
$ws.Range("I2").Value = "3V3"
$ws.Range("I2").HorizontalAlignment = -4108

# --- Row 7: DATA -> DOUT ---
$ws.Range("A7").Value = "DOUT"
$ws.Range("A7").HorizontalAlignment = -4108

# --- Remove row 19 ("Onbekend"), shifting everything below up by one ---
$ws.Rows("19").Delete()

# --- New rows 25-34 (after the shift, "SD card" block now sits on rows 21-24) ---
$ws.Range("A25").Value = "RFID"
$ws.Range("B25").Value = "SCLK"
$ws.Range("A26").Value = "RFID"
$ws.Range("B26").Value = "MISO"
$ws.Range("A27").Value = "RFID"
$ws.Range("B27").Value = "MOSI"
$ws.Range("A28").Value = "RFID"
$ws.Range("B28").Value = "CS"

$ws.Range("A29").Value = "LED rings"
$ws.Range("B29").Value = "DOUT"

$ws.Range("A30").Value = "LED display"
$ws.Range("B30").Value = "DOUT"

$ws.Range("A31").Value = "Touch"
$ws.Range("B31").Value = "CLK"
$ws.Range("A32").Value = "Touch"
$ws.Range("B32").Value = "SH/LDn"
$ws.Range("A33").Value = "Touch"
$ws.Range("B33").Value = "DIN"

$ws.Range("A34").Value = "Menutoets"
$ws.Range("B34").Value = "DIN"

$ws.Range("A25:B34").HorizontalAlignment = -4108

# Column A needs to widen to fit the new longer labels ("LED display", "Menutoets", ...)
$ws.Columns("A").AutoFit()

# Restore the cursor to A8, matching the saved selection in the workbook
$ws.Range("A8").Select()
